# Auto-generated edit script: applies updated cell values to the cryptos list (sheet1)
# per the commit diff (prices/volumes refreshed, EnergySwap/FraxShare rows swapped).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.623.23"
$ws.Range("D3").Value = "2.291.80"
$ws.Range("E3").Value = "  -1.31%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "300.19"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.46%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "96.27"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.73%  "
$ws.Range("E7").Value = "  -1.60%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  -2.85%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "33.50"
$ws.Range("D10").Style = "Normal"
$ws.Range("E11").Value = "  -0.40%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "48.81"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.86%  "
$ws.Range("E13").Value = "  +1.94%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "16.75"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +7.27%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.75"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.46%  "
$ws.Range("D16").Value = "2.649.49"
$ws.Range("E16").Value = "  -1.28%  "
$ws.Range("D17").Value = "2.307.14"
$ws.Range("E17").Value = "  -1.38%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.794"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.04%  "
$ws.Range("D19").Value = "42.589.05"
$ws.Range("E20").Value = "  -0.76%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.53"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.10%  "
$ws.Range("E22").Value = "  -1.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.81"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.98%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "236.22"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.39%  "
$ws.Range("E25").Value = "  +0.14%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.47"
$ws.Range("D26").Style = "Normal"
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "24.45"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.19"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.85%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "166.68"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.80%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.08"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.34%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.15"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.49%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.78"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.11%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.96"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.80%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.22"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.76%  "
$ws.Range("E37").Value = "  -1.45%  "
$ws.Range("E38").Value = "  -1.45%  "
$ws.Range("E39").Value = "  -3.25%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0995"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.15%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.74"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.92%  "
$ws.Range("E42").Value = "  -1.58%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.38"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -7.55%  "
$ws.Range("D44").Value = "1.957.64"
$ws.Range("E44").Value = "  -0.71%  "
$ws.Range("E45").Value = "  -1.47%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.79"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.20%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "17.52"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.66%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.82"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.64%  "
$ws.Range("D49").Value = "2.518.23"
$ws.Range("E49").Value = "  -1.11%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "52.75"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.48%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.49"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.92%  "
